$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").WrapText = $true
Write-Host "A1 style ok"
$ws.Range("A10").Copy()
$ws.Range("A370").PasteSpecial(-4122)
